# ============================================================================
# [ADDITIONAL SCRAPING] add a "Player Info" sheet (ID/NAME/BATTING_HAND/
# BOWL_STYLE for player 6783) ahead of the existing "ODI Batting" / "ODI
# Bowling" sheets, and replace MATCH_CARD_LINK (full scorecard URL) with a
# plain MATCH_CODE number on both of the existing sheets.
# ============================================================================

$wb = $excel.ActiveWorkbook

# Helper: write a value as TEXT into a cell. Bare numeric-looking strings
# (e.g. "1", "9.0", "4614") get auto-coerced to real numbers by the normal
# Value setter, which would both change their stored type and (for "9.0")
# lose the trailing zero -- so force "text entry" the same way a user would
# by typing a leading apostrophe, which Excel strips from the stored value.
function Set-CellText($range, [string]$text) {
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

$nbsp = [char]0x00A0

# ----------------------------------------------------------------------
# Step 1: grab the two existing sheets BEFORE renaming anything, so the
# object references stay valid as names/content change underneath them.
#   sheetOne currently holds "ODI Batting" data   (internal sheetId 1)
#   sheetTwo currently holds "ODI Bowling" data   (internal sheetId 2)
# Re-using these two in place (instead of deleting+re-adding) keeps their
# sheetId == 1 / 2, and a freshly Add()-ed sheet naturally becomes 3 --
# matching Player Info=1, ODI Batting=2, ODI Bowling=3 in the target file.
# ----------------------------------------------------------------------
$sheetOne = $wb.Worksheets.Item(1)
$sheetTwo = $wb.Worksheets.Item(2)

# New trailing sheet for the (new) "ODI Bowling" content.
$sheetThree = $wb.Worksheets.Add($null, $sheetTwo)

# ----------------------------------------------------------------------
# Step 2: sheetOne becomes "Player Info"
# ----------------------------------------------------------------------
$sheetOne.UsedRange.Clear()
$sheetOne.Name = "Player Info"

Set-CellText $sheetOne.Range("A1") "ID"
Set-CellText $sheetOne.Range("B1") "NAME"
Set-CellText $sheetOne.Range("C1") "BATTING_HAND"
Set-CellText $sheetOne.Range("D1") "BOWL_STYLE"

$piHeader = $sheetOne.Range("A1:D1")
$piHeader.Font.Bold = $true
$piHeader.HorizontalAlignment = -4108   # xlCenter
$piHeader.VerticalAlignment = -4160     # xlTop
$piHeader.Borders.LineStyle = 1         # xlContinuous

Set-CellText $sheetOne.Range("A2") "6783"
Set-CellText $sheetOne.Range("B2") "Graham Ian Hume"
Set-CellText $sheetOne.Range("C2") "Left Handed"
Set-CellText $sheetOne.Range("D2") "Right Arm Fast Medium"

# ----------------------------------------------------------------------
# Step 3: sheetTwo becomes "ODI Batting" (the data that used to live on
# sheetOne), with MATCH_CARD_LINK -> MATCH_CODE (bare match-code number).
# ----------------------------------------------------------------------
$sheetTwo.UsedRange.Clear()
$sheetTwo.Name = "ODI Batting"

$battingHeaders = @("MATCH_NUMBER","INNING_NUMBER","MATCH_DATE","MATCH_CODE","MATCH_INNING","OPPONENT","VENUE","DISMISSAL","RUNS_SCORED","BALLS_FACED")
$battingCols = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $battingHeaders.Count; $i++) {
    Set-CellText $sheetTwo.Range($battingCols[$i] + "1") $battingHeaders[$i]
}
$batHeaderRange = $sheetTwo.Range("A1:J1")
$batHeaderRange.Font.Bold = $true
$batHeaderRange.HorizontalAlignment = -4108
$batHeaderRange.VerticalAlignment = -4160
$batHeaderRange.Borders.LineStyle = 1

$battingRows = @(
    @("1","1","15/07/2022","4614","2nd","New Zealand","The Village","not out","7*","11"),
    @("2","","18/01/2023","4693","1st","Zimbabwe","Harare Sports Club","did not bat","-","-"),
    @("3","","21/01/2023","4694","1st","Zimbabwe","Harare Sports Club","did not bat","-","-"),
    @("4","","23/01/2023","4696",$nbsp,"Zimbabwe","Harare Sports Club","did not bat","-","-"),
    @("5","2","18/03/2023","4726","2nd","Bangladesh","Sylhet Stadium","not out","2*","5"),
    @("6","","20/03/2023","4729",$nbsp,"Bangladesh","Sylhet Stadium","did not bat","-","-"),
    @("7","3","23/03/2023","4734","1st","Bangladesh","Sylhet Stadium","lbw b Hasan Mahmud","3","19")
)

for ($r = 0; $r -lt $battingRows.Count; $r++) {
    $rowNum = $r + 2
    $rowVals = $battingRows[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $val = $rowVals[$c]
        if ($val -eq "") { continue }   # leave blank, same as source workbook
        Set-CellText $sheetTwo.Range($battingCols[$c] + $rowNum) $val
    }
}

# ----------------------------------------------------------------------
# Step 4: sheetThree becomes "ODI Bowling" (the data that used to live on
# sheetTwo), with MATCH_CARD_LINK -> MATCH_CODE (bare match-code number).
# ----------------------------------------------------------------------
$sheetThree.Name = "ODI Bowling"

$bowlingHeaders = @("MATCH_NUMBER","MATCH_CODE","MATCH_INNING","OPPONENT","VENUE","OVERS","WICKETS_RUNS")
$bowlingCols = @("A","B","C","D","E","F","G")
for ($i = 0; $i -lt $bowlingHeaders.Count; $i++) {
    Set-CellText $sheetThree.Range($bowlingCols[$i] + "1") $bowlingHeaders[$i]
}
$bowlHeaderRange = $sheetThree.Range("A1:G1")
$bowlHeaderRange.Font.Bold = $true
$bowlHeaderRange.HorizontalAlignment = -4108
$bowlHeaderRange.VerticalAlignment = -4160
$bowlHeaderRange.Borders.LineStyle = 1

$bowlingRows = @(
    @("1","4614","1st","New Zealand","The Village","9.0","0/62"),
    @("2","4693","2nd","Zimbabwe","Harare Sports Club","6.0","2/41"),
    @("3","4694","2nd","Zimbabwe","Harare Sports Club","9.0","2/52"),
    @("4","4696","1st","Zimbabwe","Harare Sports Club","5.0","0/20"),
    @("5","4726","1st","Bangladesh","Sylhet Stadium","10.0","4/60"),
    @("6","4729","1st","Bangladesh","Sylhet Stadium","10.0","3/58"),
    @("7","4734","2nd","Bangladesh","Sylhet Stadium","3.0","0/15")
)

for ($r = 0; $r -lt $bowlingRows.Count; $r++) {
    $rowNum = $r + 2
    $rowVals = $bowlingRows[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        Set-CellText $sheetThree.Range($bowlingCols[$c] + $rowNum) $rowVals[$c]
    }
}

Write-Host "Sheets now: $([string]::Join(', ', ($wb.Worksheets | ForEach-Object { $_.Name })))"
